$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 corresponds to "Tháng 8" record - update last_edited_time and the
# several numeric figures that changed as part of the new payroll report.

# last_edited_time (column D) holds a shared text value that is reused by
# several other rows (D3, D4, D5, D7, D13) pointing at the exact same
# string "2024-08-09T19:28:00.000Z". The source edit updated that shared
# string in place, so every row referencing it must be updated too in
# order to end up with the same shared string content.
$oldLastEdited = "2024-08-09T19:28:00.000Z"
$newLastEdited = "2024-08-12T02:00:00.000Z"
foreach ($r in 2..13) {
    $cell = $ws.Range("D$r")
    if ($cell.Value2 -eq $oldLastEdited) {
        $cell.Value = $newLastEdited
    }
}

# properties.Chi tiêu.number
$ws.Range("W10").Value = 11157000

# properties.Lũy kế.formula.number
$ws.Range("AA10").Value = 6843000

# properties.Tổng doanh thu.formula.number
$ws.Range("AE10").Value = 18000000

# properties.Đã thanh toán.number
$ws.Range("AH10").Value = 18000000

# properties.Số lượng đơn.number
$ws.Range("AK10").Value = 5

# properties.Đơn giá.number
$ws.Range("AQ10").Value = 19000000
